# Auto-generated edit script applying the Rafflesia_Profits.xlsx diff
# Updates cached LeveProfit/price-lookup values on the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets
# (scheduled market-data refresh - values only, no structural changes).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 78.066666
$ws.Range("I9").Value = 80.5
$ws.Range("J9").Value = 44
$ws.Range("K9").Value = 80.5
$ws.Range("L9").Value = 44
$ws.Range("M9").Value = 88.5
$ws.Range("N9").Value = -382
# Row 28
$ws.Range("H28").Value = 881.75
$ws.Range("I28").Value = 518.3333
$ws.Range("K28").Value = 518.3333
$ws.Range("M28").Value = -33.33330000000001
# Row 31
$ws.Range("H31").Value = 4311.2
$ws.Range("J31").Value = 9300
$ws.Range("L31").Value = 27900
$ws.Range("N31").Value = -28360
# Row 46
$ws.Range("H46").Value = 27300
$ws.Range("J46").Value = 24600
$ws.Range("L46").Value = 73800
$ws.Range("N46").Value = -74038
# Row 60
$ws.Range("H60").Value = 27300
$ws.Range("J60").Value = 24600
$ws.Range("L60").Value = 73800
$ws.Range("N60").Value = -74768
# Row 62
$ws.Range("H62").Value = 800
$ws.Range("I62").Value = 800
$ws.Range("K62").Value = 800
$ws.Range("M62").Value = -176
# Row 65
$ws.Range("H65").Value = 800
$ws.Range("I65").Value = 800
$ws.Range("K65").Value = 4000
$ws.Range("M65").Value = -880
# Row 74
$ws.Range("J74").Value = 2000
$ws.Range("L74").Value = 2000
$ws.Range("N74").Value = -3872
# Row 77
$ws.Range("J77").Value = 2000
$ws.Range("L77").Value = 10000
$ws.Range("N77").Value = -19360
# Row 86
$ws.Range("H86").Value = 3001
$ws.Range("I86").Value = 3003
$ws.Range("J86").Value = 2999
$ws.Range("K86").Value = 3003
$ws.Range("L86").Value = 2999
$ws.Range("M86").Value = -1880
$ws.Range("N86").Value = -5245
# Row 89
$ws.Range("H89").Value = 3001
$ws.Range("I89").Value = 3003
$ws.Range("J89").Value = 2999
$ws.Range("K89").Value = 15015
$ws.Range("L89").Value = 14995
$ws.Range("M89").Value = -9399
$ws.Range("N89").Value = -26227
# Row 98
$ws.Range("H98").Value = 896
$ws.Range("I98").Value = 896
$ws.Range("K98").Value = 896
$ws.Range("M98").Value = 602
# Row 122
$ws.Range("H122").Value = 896
$ws.Range("I122").Value = 896
$ws.Range("K122").Value = 2688
$ws.Range("M122").Value = -238
# Row 135
$ws.Range("H135").Value = 2099.6667
$ws.Range("I135").Value = 2099.6667
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 18897.0003
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -16362.0003

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2249.75
$ws.Range("I45").Value = 1500
$ws.Range("K45").Value = 1500
$ws.Range("M45").Value = -1123
# Row 122
$ws.Range("H122").Value = 3936.4285
$ws.Range("I122").Value = 1166.3334
$ws.Range("K122").Value = 3499.0002
$ws.Range("M122").Value = -1049.0002
# Row 132
$ws.Range("H132").Value = 1300
$ws.Range("I132").Value = 1300
$ws.Range("K132").Value = 3900
$ws.Range("M132").Value = -1370

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 2997.5
$ws.Range("I80").Value = 2997.5
$ws.Range("K80").Value = 2997.5
$ws.Range("M80").Value = -1999.5
# Row 83
$ws.Range("H83").Value = 2997.5
$ws.Range("I83").Value = 2997.5
$ws.Range("K83").Value = 14987.5
$ws.Range("M83").Value = -9995.5
# Row 99
$ws.Range("H99").Value = 1911.125
$ws.Range("I99").Value = 1898.4286
$ws.Range("K99").Value = 1898.4286
$ws.Range("M99").Value = -400.4286
# Row 105
$ws.Range("H105").Value = 13267.929
$ws.Range("I105").Value = 13267.929
$ws.Range("K105").Value = 13267.929
$ws.Range("M105").Value = -11520.929
# Row 107
$ws.Range("H107").Value = 1132.125
$ws.Range("I107").Value = 416.2
$ws.Range("K107").Value = 416.2
$ws.Range("M107").Value = 1503.8

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("N16").Value = 0
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").ClearContents()
$ws.Range("N57").Value = 0
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("N113").Value = 0
# Row 134
$ws.Range("H134").Value = 5849.4
$ws.Range("I134").Value = 3082.6667
$ws.Range("K134").Value = 9248.000100000001
$ws.Range("M134").Value = -6713.000100000001

$ws = $wb.Worksheets.Item("GSM")
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
# Row 58
$ws.Range("H58").Value = 2422.5
$ws.Range("I58").Value = 2422.5
$ws.Range("K58").Value = 2422.5
$ws.Range("M58").Value = -2145.5
# Row 63
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372
# Row 66
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864
# Row 97
$ws.Range("H97").Value = 636.875
$ws.Range("I97").Value = 636.875
$ws.Range("K97").Value = 636.875
$ws.Range("M97").Value = -140.875
# Row 113
$ws.Range("H113").Value = 3093.889
$ws.Range("I113").Value = 2190.2856
$ws.Range("J113").Value = 6256.5
$ws.Range("K113").Value = 2190.2856
$ws.Range("L113").Value = 6256.5
$ws.Range("M113").Value = -20.28560000000016
$ws.Range("N113").Value = -10596.5

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 3942.7144
$ws.Range("I16").Value = 3942.7144
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3942.7144
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -3772.7144
# Row 22
$ws.Range("H22").Value = 3375
$ws.Range("I22").Value = 1750
$ws.Range("K22").Value = 1750
$ws.Range("M22").Value = -1455
# Row 27
$ws.Range("H27").Value = 3375
$ws.Range("I27").Value = 1750
$ws.Range("K27").Value = 1750
$ws.Range("M27").Value = -1643
# Row 93
$ws.Range("H93").Value = 2842
$ws.Range("I93").Value = 777
$ws.Range("J93").Value = 3874.5
$ws.Range("K93").Value = 777
$ws.Range("L93").Value = 3874.5
$ws.Range("M93").Value = 471
$ws.Range("N93").Value = -6370.5
# Row 100
$ws.Range("H100").Value = 1956
$ws.Range("I100").Value = 1956
$ws.Range("K100").Value = 1956
$ws.Range("M100").Value = -1415

$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Range("H63").Value = 10000
$ws.Range("I63").Value = 10000
$ws.Range("K63").Value = 10000
$ws.Range("M63").Value = -9376
# Row 66
$ws.Range("H66").Value = 10000
$ws.Range("I66").Value = 10000
$ws.Range("K66").Value = 30000
$ws.Range("M66").Value = -26880
# Row 81
$ws.Range("H81").Value = 400
$ws.Range("I81").Value = 400
$ws.Range("K81").Value = 800
$ws.Range("M81").Value = 261
# Row 84
$ws.Range("H84").Value = 400
$ws.Range("I84").Value = 400
$ws.Range("K84").Value = 4000
$ws.Range("M84").Value = 1304
# Row 107
$ws.Range("H107").Value = 4042.3333
$ws.Range("I107").Value = 1622.5
$ws.Range("J107").Value = 5252.25
$ws.Range("K107").Value = 4867.5
$ws.Range("L107").Value = 15756.75
$ws.Range("M107").Value = -2947.5
$ws.Range("N107").Value = -19596.75

Write-Output "Updated 184 cells, cleared 6 cells across 7 sheets."